# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (all rows refer to this sheet's own row numbering)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 1219
$wsExpo.Range("F4").Value  = 14620
$wsExpo.Range("F5").Value  = 17747
$wsExpo.Range("F16").Value = 48
$wsExpo.Range("F17").Value = 155
$wsExpo.Range("F18").Value = 43
$wsExpo.Range("F24").Value = 7314
$wsExpo.Range("F30").Value = 5871
$wsExpo.Range("F31").Value = 68
$wsExpo.Range("F33").Value = 143

# Sheet "全部类型" (same events, but row numbers are offset by extra rows
# present only on this combined sheet)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 1219
$wsAll.Range("F4").Value  = 14620
$wsAll.Range("F5").Value  = 17747
$wsAll.Range("F16").Value = 48
$wsAll.Range("F17").Value = 155
$wsAll.Range("F18").Value = 43
$wsAll.Range("F25").Value = 7314
$wsAll.Range("F32").Value = 5871
$wsAll.Range("F33").Value = 68
$wsAll.Range("F35").Value = 143
